$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell F1, matching the style of the other header cells (E1).
$ws.Range("F1").Value = "time_taken"
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$timestamps = @(
    "2021-10-05 13:42:01.330939",
    "2021-10-05 13:42:01.330950",
    "2021-10-05 13:42:01.330954",
    "2021-10-05 13:42:01.330956",
    "2021-10-05 13:42:01.330959",
    "2021-10-05 13:42:01.330962",
    "2021-10-05 13:42:01.330965",
    "2021-10-05 13:42:01.330967",
    "2021-10-05 13:42:01.330970",
    "2021-10-05 13:42:01.330973",
    "2021-10-05 13:42:01.330975",
    "2021-10-05 13:42:01.330978",
    "2021-10-05 13:42:01.330980",
    "2021-10-05 13:42:01.330983",
    "2021-10-05 13:42:01.330986",
    "2021-10-05 13:42:01.330988",
    "2021-10-05 13:42:01.330991",
    "2021-10-05 13:42:01.330993",
    "2021-10-05 13:42:01.330996",
    "2021-10-05 13:42:01.330999",
    "2021-10-05 13:42:01.331001",
    "2021-10-05 13:42:01.331004",
    "2021-10-05 13:42:01.331006",
    "2021-10-05 13:42:01.331009",
    "2021-10-05 13:42:01.331012",
    "2021-10-05 13:42:01.331014",
    "2021-10-05 13:42:01.331017",
    "2021-10-05 13:42:01.331019",
    "2021-10-05 13:42:01.331022",
    "2021-10-05 13:42:01.331025",
    "2021-10-05 13:42:01.331027",
    "2021-10-05 13:42:01.331030",
    "2021-10-05 13:42:01.331033",
    "2021-10-05 13:42:01.331035",
    "2021-10-05 13:42:01.331038",
    "2021-10-05 13:42:01.331040",
    "2021-10-05 13:42:01.331043",
    "2021-10-05 13:42:01.331046",
    "2021-10-05 13:42:01.331048",
    "2021-10-05 13:42:01.331051",
    "2021-10-05 13:42:01.331054",
    "2021-10-05 13:42:01.331057",
    "2021-10-05 13:42:01.331059",
    "2021-10-05 13:42:01.331062",
    "2021-10-05 13:42:01.331064",
    "2021-10-05 13:42:01.331067",
    "2021-10-05 13:42:01.331069",
    "2021-10-05 13:42:01.331072",
    "2021-10-05 13:42:01.331074",
    "2021-10-05 13:42:01.331077",
    "2021-10-05 13:42:01.331080",
    "2021-10-05 13:42:01.331082",
    "2021-10-05 13:42:01.331085",
    "2021-10-05 13:42:01.331088",
    "2021-10-05 13:42:01.331091",
    "2021-10-05 13:42:01.331093",
    "2021-10-05 13:42:01.331096",
    "2021-10-05 13:42:01.331098",
    "2021-10-05 13:42:01.331101",
    "2021-10-05 13:42:01.331103",
    "2021-10-05 13:42:01.331106"
)

for ($i = 0; $i -lt $timestamps.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $timestamps[$i]
}
